$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks (Vahan/Eurotouram/Lilit/Levon addresses) so we
# can lay down the new recipient list cleanly.
$ws.Hyperlinks.Delete()

# Drop the three trailing recipient rows (lil-3@mail.ru, lilis88@mail.ru,
# levon.eurotour@gmail.com) - only three recipients remain after the fix.
$ws.Rows("5:7").Delete()

# Replace the remaining recipients with the corrected mailing list.
$ws.Range("A2").Value = "gevorgadamyan@yahoo.com"
$ws.Range("A3").Value = "adamyangevorg4@gmail.com"
$ws.Range("A4").Value = "gevorgadamyan@outlook.com"

# Re-create the mailto hyperlinks for the new addresses.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:gevorgadamyan@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:adamyangevorg4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:gevorgadamyan@outlook.com")

# Adding the hyperlinks re-applies Excel's default "Hyperlink" look; restore
# the original left-aligned formatting that the sheet already used.
$ws.Range("A2:A4").HorizontalAlignment = -4131

$ws.Range("G8").Select() | Out-Null
